$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Count" header into the new E1 cell (F column already has it)
$ws.Range("E1").Value = "Count"

# The DDS accumulator clock reference was previously drifting per-row
# (50000001..50000007); the "baud rate switch" fix pins every row back
# to the master 50MHz clock so the computed tuning-word counts line up.
$ws.Range("C4").Value = 50000000
$ws.Range("C5").Value = 50000000
$ws.Range("C6").Value = 50000000
$ws.Range("C7").Value = 50000000
$ws.Range("C8").Value = 50000000
$ws.Range("C9").Value = 50000000
$ws.Range("C10").Value = 50000000

# Reselect F10 (the last "Count" cell) and zoom in, as in the saved view
$ws.Range("F10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 390
